$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update cell values per the diff
$ws.Range("R8").Value = 1
$ws.Range("R12").Value = 1

$ws.Range("L13").Value = 1
$ws.Range("O13").Value = 1
$ws.Range("P13").Value = 1
$ws.Range("Q13").Value = 1
$ws.Range("S13").Value = 1

$ws.Range("L17").Value = 1
$ws.Range("N17").Value = 1
$ws.Range("O17").Value = 1
$ws.Range("P17").Value = 1
$ws.Range("Q17").Value = 1
$ws.Range("R17").Value = 1
$ws.Range("S17").Value = 1

# Update the view: scroll so column G is the leftmost visible column,
# and update the active selection to R40.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("R40").Select()
